$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.1815
$ws.Range("B4").Value = 4.292500000000003
$ws.Range("D4").Value = -7.042000000000004

$ws.Range("B5").Value = 4.9883

$ws.Range("A7").Value = -21.55380000000001

$ws.Range("B8").Value = 4.943199999999999

$ws.Range("D9").Value = -8.348600000000006

$ws.Range("A16").Value = -21.45540000000002
$ws.Range("B16").Value = 4.834300000000002

$ws.Range("D18").Value = -8.366699999999993
